$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Arun's email (B4) has become an exception - clear the value but keep the
# cell's existing formatting, and remove the mailto: hyperlink that was
# attached to it.
$ws.Range("B4").ClearContents()

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$4') {
        $hl.Delete()
    }
}

# Leave the selection where the author left it when they saved.
$ws.Range("C21").Select() | Out-Null
